$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrected parameter values ---
$ws.Range("C21").Value = 5        # c - center to bearing
$ws.Range("C22").Value = 0.3      # D[S] - shaft-diameter
$ws.Range("C32").Value = 3        # omega[0] - eigenfrequ. of blade
$ws.Range("C34").Value = 0.1      # A[B] - blade cross-sec area
$ws.Range("C36").Value = 0.1      # N[0]

# --- Update view / selection to match where the author left off ---
$ws.Range("C35").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
